# Rename database-field labels:
#   "NR_USUA"    -> "CD_FUNC"     (cadr_func!A2, cadr_empr_parc!A7)
#   "FK_NR_USUA" -> "FK_CD_USUA"  (sens_info!A10)
# and refresh the saved view state (selected ranges / active sheet) to match
# what Excel records after the user clicked through the sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("cadr_func")
$ws2 = $wb.Worksheets.Item("sens_info")
$ws3 = $wb.Worksheets.Item("cadr_empr_parc")

# --- data edits -----------------------------------------------------------
$ws1.Range("A2").Value = "CD_FUNC"
$ws3.Range("A7").Value = "CD_FUNC"
$ws2.Range("A10").Value = "FK_CD_USUA"

# --- view state -------------------------------------------------------------
# Select the full used range on each sheet (as Excel does on Ctrl+A), then
# finish with "sens_info" as the active / selected tab.
$ws1.Activate()
$ws1.Range("A1:E5").Select()

$ws3.Activate()
$ws3.Range("A1:E7").Select()

$ws2.Activate()
$ws2.Range("A1:E11").Select()

$wb.Save()
